# Insert a new slide "So what about data?" right before the final
# "Let's do this" slide (currently the last slide, #22), reusing the
# same "Title and Content" layout.

$p = $ppt.ActivePresentation

$lastSlide = $p.Slides.Item($p.Slides.Count)
$layout = $lastSlide.CustomLayout

# Insert the new slide right before the current last slide.
$newSlide = $p.Slides.AddSlide($lastSlide.SlideIndex, $layout)

# Title placeholder.
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "So what about data?"

# Body / content placeholder - three text paragraphs separated by blank lines.
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Git is not a data platform" + "`r" + "`r" + `
    "But used to version control workflows surrounding data" + "`r" + "`r" + `
    "Git workflow stimulates documentation, meta data, accessibility, transparency"
